$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 627.30304
$ws.Range("I98").Value = 627.30304
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 627.30304
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 870.69696
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 627.30304
$ws.Range("I122").Value = 627.30304
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1881.90912
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 568.09088
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1659.28
$ws.Range("I132").Value = 1140.5862
$ws.Range("J132").Value = 3428.9412
$ws.Range("K132").Value = 3421.7586
$ws.Range("L132").Value = 10286.8236
$ws.Range("M132").Value = -891.7586000000001
$ws.Range("N132").Value = -15346.8236

$ws.Range("H138").Value = 2365.3704
$ws.Range("I138").Value = 997.875
$ws.Range("J138").Value = 4354.4546
$ws.Range("K138").Value = 2993.625
$ws.Range("L138").Value = 13063.3638
$ws.Range("M138").Value = 2146.375
$ws.Range("N138").Value = -23343.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4021.9456
$ws.Range("I32").Value = 4119.365
$ws.Range("K32").Value = 4119.365
$ws.Range("M32").Value = -3832.365

$ws.Range("H61").Value = 1727.1852
$ws.Range("I61").Value = 1314.4117
$ws.Range("K61").Value = 1314.4117
$ws.Range("M61").Value = -1102.4117

$ws.Range("H122").Value = 1548.1875
$ws.Range("I122").Value = 1471.1111
$ws.Range("J122").Value = 1647.2858
$ws.Range("K122").Value = 4413.3333
$ws.Range("L122").Value = 4941.857400000001
$ws.Range("M122").Value = -1963.3333
$ws.Range("N122").Value = -9841.857400000001

$ws.Range("H136").Value = 1727.1852
$ws.Range("I136").Value = 1314.4117
$ws.Range("K136").Value = 3943.2351
$ws.Range("M136").Value = -1393.2351

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 30000.375
$ws.Range("J63").Value = 30000.375
$ws.Range("L63").Value = 30000.375
$ws.Range("N63").Value = -31372.375

$ws.Range("H66").Value = 30000.375
$ws.Range("J66").Value = 30000.375
$ws.Range("L66").Value = 90001.125
$ws.Range("N66").Value = -96865.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17144.338
$ws.Range("I31").Value = 1021.2683
$ws.Range("J31").Value = 39179.2
$ws.Range("K31").Value = 1021.2683
$ws.Range("L31").Value = 39179.2
$ws.Range("M31").Value = -726.2683
$ws.Range("N31").Value = -39769.2

$ws.Range("H34").Value = 17144.338
$ws.Range("I34").Value = 1021.2683
$ws.Range("J34").Value = 39179.2
$ws.Range("K34").Value = 1021.2683
$ws.Range("L34").Value = 39179.2
$ws.Range("M34").Value = -819.2683
$ws.Range("N34").Value = -39583.2

$ws.Range("H86").Value = 3679.8462
$ws.Range("I86").Value = 3338.4614
$ws.Range("J86").Value = 4021.2307
$ws.Range("K86").Value = 3338.4614
$ws.Range("L86").Value = 4021.2307
$ws.Range("M86").Value = -2215.4614
$ws.Range("N86").Value = -6267.2307

$ws.Range("H89").Value = 3679.8462
$ws.Range("I89").Value = 3338.4614
$ws.Range("J89").Value = 4021.2307
$ws.Range("K89").Value = 16692.307
$ws.Range("L89").Value = 20106.1535
$ws.Range("M89").Value = -11076.307
$ws.Range("N89").Value = -31338.1535

$ws.Range("H122").Value = 886.36
$ws.Range("I122").Value = 897.4375
$ws.Range("J122").Value = 866.6667
$ws.Range("K122").Value = 2692.3125
$ws.Range("L122").Value = 2600.0001
$ws.Range("M122").Value = -242.3125
$ws.Range("N122").Value = -7500.0001

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 694.7439000000001
$ws.Range("I131").Value = 458.9565
$ws.Range("J131").Value = 996.0278
$ws.Range("K131").Value = 1376.8695
$ws.Range("L131").Value = 2988.0834
$ws.Range("M131").Value = 3663.1305
$ws.Range("N131").Value = -13068.0834

$ws.Range("H132").Value = 994.2308
$ws.Range("I132").Value = 941.5
$ws.Range("K132").Value = 8473.5
$ws.Range("M132").Value = -5943.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1200.75
$ws.Range("I122").Value = 1201
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3603
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1153
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 700.6667
$ws.Range("I16").Value = 700
$ws.Range("J16").Value = 701
$ws.Range("K16").Value = 700
$ws.Range("L16").Value = 701
$ws.Range("M16").Value = -530
$ws.Range("N16").Value = -1041

$ws.Range("H40").Value = 2141.3572
$ws.Range("I40").Value = 1899.8572
$ws.Range("J40").Value = 2382.8572
$ws.Range("K40").Value = 1899.8572
$ws.Range("L40").Value = 2382.8572
$ws.Range("M40").Value = -1763.8572
$ws.Range("N40").Value = -2654.8572

$ws.Range("H58").Value = 7000
$ws.Range("I58").Value = 7000
$ws.Range("J58").Value = 7000
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = -6740
$ws.Range("N58").Value = -7520

$ws.Range("H122").Value = 42835.28
$ws.Range("I122").Value = 69164.13
$ws.Range("J122").Value = 3342
$ws.Range("K122").Value = 207492.39
$ws.Range("L122").Value = 10026
$ws.Range("M122").Value = -205042.39
$ws.Range("N122").Value = -14926

$ws.Range("H132").Value = 4196.5137
$ws.Range("I132").Value = 4339.7812
$ws.Range("J132").Value = 3279.6
$ws.Range("K132").Value = 13019.3436
$ws.Range("L132").Value = 9838.799999999999
$ws.Range("M132").Value = -10489.3436
$ws.Range("N132").Value = -14898.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1073.3334
$ws.Range("I81").Value = 1041.4286
$ws.Range("J81").Value = 1137.1428
$ws.Range("K81").Value = 2082.8572
$ws.Range("L81").Value = 2274.2856
$ws.Range("M81").Value = -1021.8572
$ws.Range("N81").Value = -4396.2856

$ws.Range("H84").Value = 1073.3334
$ws.Range("I84").Value = 1041.4286
$ws.Range("J84").Value = 1137.1428
$ws.Range("K84").Value = 10414.286
$ws.Range("L84").Value = 11371.428
$ws.Range("M84").Value = -5110.286
$ws.Range("N84").Value = -21979.428

$ws.Range("H126").Value = 1333.5714
$ws.Range("I126").Value = 1050
$ws.Range("J126").Value = 1546.25
$ws.Range("K126").Value = 3150
$ws.Range("L126").Value = 4638.75
$ws.Range("M126").Value = -680
$ws.Range("N126").Value = -9578.75
